$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.359656
$ws.Range("H2").Value = 7.078968
$ws.Range("I2").Value = 0.1135804410355361
$ws.Range("J2").Value = 0.1135804410355361
$ws.Range("O2").Value = 0.9349445792302935
$ws.Range("P2").Value = 0.9349445792302935
$ws.Range("Q2").Value = 1.545847613544
$ws.Range("R2").Value = 13.912628521896
$ws.Range("S2").Value = 0.1061914176527605
$ws.Range("T2").Value = 0.1061914176527605
$ws.Range("G3").Value = 2.359656
$ws.Range("H3").Value = 7.078968
$ws.Range("I3").Value = 0.1135804410355361
$ws.Range("J3").Value = 0.1135804410355361
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.04558433333333334
$ws.Range("N3").Value = 0.136753
$ws.Range("O3").Value = 0.0650554207697065
$ws.Range("P3").Value = 0.06505542076970648
$ws.Range("Q3").Value = 0.107563345656
$ws.Range("R3").Value = 0.9680701109040001
$ws.Range("S3").Value = 0.007389023382775641
$ws.Range("T3").Value = 0.007389023382775639
$ws.Range("I4").Value = 0.07630393871923234
$ws.Range("J4").Value = 0.07630393871923234
$ws.Range("O4").Value = 0.9349445792302935
$ws.Range("P4").Value = 0.9349445792302935
$ws.Range("S4").Value = 0.07133995387946679
$ws.Range("T4").Value = 0.07133995387946679
$ws.Range("I5").Value = 0.07630393871923234
$ws.Range("J5").Value = 0.07630393871923234
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04558433333333334
$ws.Range("N5").Value = 0.136753
$ws.Range("O5").Value = 0.0650554207697065
$ws.Range("P5").Value = 0.06505542076970648
$ws.Range("Q5").Value = 0.07226162234044445
$ws.Range("R5").Value = 0.6503546010640001
$ws.Range("S5").Value = 0.004963984839765559
$ws.Range("T5").Value = 0.004963984839765559
$ws.Range("G6").Value = 3.018243333333333
$ws.Range("H6").Value = 9.054729999999999
$ws.Range("I6").Value = 0.1452810956141771
$ws.Range("J6").Value = 0.1452810956141771
$ws.Range("O6").Value = 0.9349445792302935
$ws.Range("P6").Value = 0.9349445792302935
$ws.Range("Q6").Value = 1.977298493478889
$ws.Range("R6").Value = 17.79568644131
$ws.Range("S6").Value = 0.1358297728091128
$ws.Range("T6").Value = 0.1358297728091128
$ws.Range("G7").Value = 3.018243333333333
$ws.Range("H7").Value = 9.054729999999999
$ws.Range("I7").Value = 0.1452810956141771
$ws.Range("J7").Value = 0.1452810956141771
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.04558433333333334
$ws.Range("N7").Value = 0.136753
$ws.Range("O7").Value = 0.0650554207697065
$ws.Range("P7").Value = 0.06505542076970648
$ws.Range("Q7").Value = 0.1375846101877778
$ws.Range("R7").Value = 1.23826149169
$ws.Range("S7").Value = 0.009451322805064251
$ws.Range("T7").Value = 0.009451322805064249
$ws.Range("G8").Value = 1.732629
$ws.Range("H8").Value = 5.197887
$ws.Range("I8").Value = 0.08339892169492499
$ws.Range("J8").Value = 0.08339892169492499
$ws.Range("O8").Value = 0.9349445792302935
$ws.Range("P8").Value = 0.9349445792302935
$ws.Range("Q8").Value = 1.135072402421
$ws.Range("R8").Value = 10.215651621789
$ws.Range("S8").Value = 0.07797336975232184
$ws.Range("T8").Value = 0.07797336975232184
$ws.Range("G9").Value = 1.732629
$ws.Range("H9").Value = 5.197887
$ws.Range("I9").Value = 0.08339892169492499
$ws.Range("J9").Value = 0.08339892169492499
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.04558433333333334
$ws.Range("N9").Value = 0.136753
$ws.Range("O9").Value = 0.0650554207697065
$ws.Range("P9").Value = 0.06505542076970648
$ws.Range("Q9").Value = 0.07898073787900001
$ws.Range("R9").Value = 0.710826640911
$ws.Range("S9").Value = 0.005425551942603149
$ws.Range("T9").Value = 0.005425551942603148
$ws.Range("G10").Value = 9.589644
$ws.Range("H10").Value = 28.768932
$ws.Range("I10").Value = 0.4615910094072114
$ws.Range("J10").Value = 0.4615910094072114
$ws.Range("O10").Value = 0.9349445792302935
$ws.Range("P10").Value = 0.9349445792302935
$ws.Range("Q10").Value = 6.282326022156
$ws.Range("R10").Value = 56.540934199404
$ws.Range("S10").Value = 0.4315620120667117
$ws.Range("T10").Value = 0.4315620120667117
$ws.Range("G11").Value = 9.589644
$ws.Range("H11").Value = 28.768932
$ws.Range("I11").Value = 0.4615910094072114
$ws.Range("J11").Value = 0.4615910094072114
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.04558433333333334
$ws.Range("N11").Value = 0.136753
$ws.Range("O11").Value = 0.0650554207697065
$ws.Range("P11").Value = 0.06505542076970648
$ws.Range("Q11").Value = 0.437137528644
$ws.Range("R11").Value = 3.934237757796
$ws.Range("S11").Value = 0.03002899734049969
$ws.Range("T11").Value = 0.03002899734049968
$ws.Range("G12").Value = 2.489795
$ws.Range("H12").Value = 7.469385
$ws.Range("I12").Value = 0.1198445935289181
$ws.Range("J12").Value = 0.1198445935289181
$ws.Range("O12").Value = 0.9349445792302935
$ws.Range("P12").Value = 0.9349445792302935
$ws.Range("Q12").Value = 1.631103711288333
$ws.Range("R12").Value = 14.679933401595
$ws.Range("S12").Value = 0.1120480530699199
$ws.Range("T12").Value = 0.1120480530699199
$ws.Range("G13").Value = 2.489795
$ws.Range("H13").Value = 7.469385
$ws.Range("I13").Value = 0.1198445935289181
$ws.Range("J13").Value = 0.1198445935289181
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.04558433333333334
$ws.Range("N13").Value = 0.136753
$ws.Range("O13").Value = 0.0650554207697065
$ws.Range("P13").Value = 0.06505542076970648
$ws.Range("Q13").Value = 0.1134956452116667
$ws.Range("R13").Value = 1.021460806905
$ws.Range("S13").Value = 0.007796540458998209
$ws.Range("T13").Value = 0.007796540458998207
